$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H36").Value = -1
$ws.Range("H37").Value = -1
$ws.Range("H38").Value = -1
$ws.Range("H42").Value = 1
$ws.Range("H45").Value = -1
$ws.Range("H49").Value = -1
